# Auto-generated edit script: updates currentAveragePrice / Leve price / profit
# columns (H-N) across rows in multiple sheets of the Unicorn_Profits workbook,
# reflecting refreshed market-board price data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 3
$ws.Range("H3").Value = 31125
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 31125
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 31125
$ws.Range("N3").Value = -31353

# ALC row 32
$ws.Range("H32").Value = 5808405
$ws.Range("I32").Value = 262.5
$ws.Range("J32").Value = 8712476
$ws.Range("K32").Value = 262.5
$ws.Range("L32").Value = 8712476
$ws.Range("M32").Value = 63.5
$ws.Range("N32").Value = -8713128

# ALC row 102
$ws.Range("H102").Value = 31125
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 31125
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 31125
$ws.Range("N102").Value = -37615

$ws = $wb.Worksheets.Item("ARM")
# ARM row 122
$ws.Range("H122").Value = 6791.6924
$ws.Range("I122").Value = 6976.8887
$ws.Range("J122").Value = 6375
$ws.Range("K122").Value = 20930.6661
$ws.Range("L122").Value = 19125
$ws.Range("M122").Value = -18480.6661
$ws.Range("N122").Value = -24025

$ws = $wb.Worksheets.Item("BSM")
# BSM row 12
$ws.Range("H12").Value = 308.77777
$ws.Range("I12").Value = 308.77777
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 308.77777
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -140.77777
$ws.Range("N12").ClearContents()

# BSM row 24
$ws.Range("H24").Value = 579
$ws.Range("I24").Value = 398.75
$ws.Range("J24").Value = 1300
$ws.Range("K24").Value = 398.75
$ws.Range("L24").Value = 1300
$ws.Range("M24").Value = -163.75
$ws.Range("N24").Value = -1770

# BSM row 25
$ws.Range("H25").Value = 1121.3334
$ws.Range("I25").Value = 382
$ws.Range("J25").Value = 2600
$ws.Range("K25").Value = 382
$ws.Range("L25").Value = 2600
$ws.Range("M25").Value = -147
$ws.Range("N25").Value = -3070

# BSM row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# BSM row 134
$ws.Range("H134").Value = 2232.9285
$ws.Range("I134").Value = 2040.5745
$ws.Range("J134").Value = 2626
$ws.Range("K134").Value = 6121.7235
$ws.Range("L134").Value = 7878
$ws.Range("M134").Value = -3586.7235
$ws.Range("N134").Value = -12948

$ws = $wb.Worksheets.Item("CRP")
# CRP row 5
$ws.Range("H5").Value = 424.46155
$ws.Range("I5").Value = 138.33333
$ws.Range("J5").Value = 669.7143
$ws.Range("K5").Value = 138.33333
$ws.Range("L5").Value = 669.7143
$ws.Range("M5").Value = -26.33332999999999
$ws.Range("N5").Value = -893.7143

# CRP row 16
$ws.Range("H16").Value = 673.2
$ws.Range("I16").Value = 632.9286
$ws.Range("J16").Value = 767.1667
$ws.Range("K16").Value = 632.9286
$ws.Range("L16").Value = 767.1667
$ws.Range("M16").Value = -345.9286
$ws.Range("N16").Value = -1341.1667

# CRP row 19
$ws.Range("H19").Value = 400122.1
$ws.Range("I19").Value = 400122.1
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 400122.1
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -399952.1
$ws.Range("N19").ClearContents()

# CRP row 22
$ws.Range("H22").Value = 631.25
$ws.Range("I22").Value = 676.6667
$ws.Range("J22").Value = 604
$ws.Range("K22").Value = 676.6667
$ws.Range("L22").Value = 604
$ws.Range("M22").Value = -326.6667
$ws.Range("N22").Value = -1304

# CRP row 24
$ws.Range("H24").Value = 400122.1
$ws.Range("I24").Value = 400122.1
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 400122.1
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -399952.1
$ws.Range("N24").ClearContents()

# CRP row 99
$ws.Range("H99").Value = 47124.59
$ws.Range("I99").Value = 101334.1
$ws.Range("J99").Value = 1950
$ws.Range("K99").Value = 101334.1
$ws.Range("L99").Value = 1950
$ws.Range("M99").Value = -99836.10000000001
$ws.Range("N99").Value = -4946

# CRP row 113
$ws.Range("H113").Value = 673.2
$ws.Range("I113").Value = 632.9286
$ws.Range("J113").Value = 767.1667
$ws.Range("K113").Value = 632.9286
$ws.Range("L113").Value = 767.1667
$ws.Range("M113").Value = 1537.0714
$ws.Range("N113").Value = -5107.1667

# CRP row 122
$ws.Range("H122").Value = 5314.143
$ws.Range("I122").Value = 10500
$ws.Range("J122").Value = 3239.8
$ws.Range("K122").Value = 31500
$ws.Range("L122").Value = 9719.400000000001
$ws.Range("M122").Value = -29050
$ws.Range("N122").Value = -14619.4

# CRP row 126
$ws.Range("H126").Value = 47124.59
$ws.Range("I126").Value = 101334.1
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 304002.3
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -301532.3
$ws.Range("N126").Value = -10790

# CRP row 132
$ws.Range("H132").Value = 1284.7703
$ws.Range("I132").Value = 744.1607
$ws.Range("J132").Value = 2966.6667
$ws.Range("K132").Value = 2232.4821
$ws.Range("L132").Value = 8900.000100000001
$ws.Range("M132").Value = 297.5178999999998
$ws.Range("N132").Value = -13960.0001

# CRP row 134
$ws.Range("H134").Value = 1394.3334
$ws.Range("I134").Value = 849.8982999999999
$ws.Range("J134").Value = 2541.5356
$ws.Range("K134").Value = 2549.6949
$ws.Range("L134").Value = 7624.6068
$ws.Range("M134").Value = -14.69489999999996
$ws.Range("N134").Value = -12694.6068

$ws = $wb.Worksheets.Item("CUL")
# CUL row 107
$ws.Range("H107").Value = 384.37143
$ws.Range("I107").Value = 228.88889
$ws.Range("J107").Value = 549
$ws.Range("K107").Value = 686.6666700000001
$ws.Range("L107").Value = 1647
$ws.Range("M107").Value = 1233.33333
$ws.Range("N107").Value = -5487

$ws = $wb.Worksheets.Item("GSM")
# GSM row 13
$ws.Range("H13").Value = 618.5
$ws.Range("I13").Value = 324.66666
$ws.Range("J13").Value = 1500
$ws.Range("K13").Value = 324.66666
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = -185.66666
$ws.Range("N13").Value = -1778

# GSM row 122
$ws.Range("H122").Value = 1320.625
$ws.Range("I122").Value = 749
$ws.Range("J122").Value = 1663.6
$ws.Range("K122").Value = 2247
$ws.Range("L122").Value = 4990.799999999999
$ws.Range("M122").Value = 203
$ws.Range("N122").Value = -9890.799999999999

# GSM row 126
$ws.Range("H126").Value = 1846.4814
$ws.Range("I126").Value = 1565
$ws.Range("J126").Value = 2650.7144
$ws.Range("K126").Value = 4695
$ws.Range("L126").Value = 7952.1432
$ws.Range("M126").Value = -2225
$ws.Range("N126").Value = -12892.1432

$ws = $wb.Worksheets.Item("LTW")
# LTW row 41
$ws.Range("H41").Value = 16912.334
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 16912.334
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 16912.334
$ws.Range("N41").Value = -17788.334

$ws = $wb.Worksheets.Item("WVR")
# WVR row 126
$ws.Range("H126").Value = 1141.826
$ws.Range("I126").Value = 814.8823
$ws.Range("J126").Value = 2068.1667
$ws.Range("K126").Value = 2444.6469
$ws.Range("L126").Value = 6204.500100000001
$ws.Range("M126").Value = 25.35310000000027
$ws.Range("N126").Value = -11144.5001

# WVR row 136
$ws.Range("H136").Value = 12799150
$ws.Range("I136").Value = 18201326
$ws.Range("J136").Value = 419162.9
$ws.Range("K136").Value = 54603978
$ws.Range("L136").Value = 1257488.7
$ws.Range("M136").Value = -54601428
$ws.Range("N136").Value = -1262588.7
